# ExportMode-byClass.xlsx - "Volt VAR WATT working"
#
# The sheet held a small 3-row export-list table (Class/Property/Column1)
# listing Loads/NumPhases/kW and PVsystems/Powers/<blank>. The commit
# collapses it to a single data row describing a PVsystems export
# (VoltagesMagAng / Powers) and swaps the Property/Column1 column order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "PVsystems | Powers |" row (row 4) entirely - the table
# shrinks from A2:C4 to A2:C3 and everything below shifts up.
$ws.Rows.Item(4).Delete() | Out-Null

# Row 1 used to carry a stray styled-but-empty cell out in column G; clear
# it so the sheet's used range shrinks back down to column F.
$ws.Range("G1").Clear() | Out-Null

# Header row: swap the Property/Column1 header labels (Column1 now sits in
# B, Property moves to C).
$ws.Range("B2").Value = "Column1"
$ws.Range("C2").Value = "Property"

# Data row: now describes the PVsystems/VoltagesMagAng/Powers export.
$ws.Range("A3").Value = "PVsystems"
$ws.Range("B3").Value = "VoltagesMagAng"
$ws.Range("C3").Value = "Powers"

# B3 previously had no explicit style (it held "NumPhases" unstyled); the
# new VoltagesMagAng value picks up the same formatting as the rest of the
# data row, so copy A3's format onto it.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Keep the table (Table1) in sync: same column swap as the header row.
$lo = $ws.ListObjects.Item(1)
$lo.ListColumns.Item(2).Name = "__swap_tmp__"
$lo.ListColumns.Item(3).Name = "Property"
$lo.ListColumns.Item(2).Name = "Column1"

# Leave the saved selection where the author left it.
$ws.Range("E13").Select() | Out-Null
